# Apply updated crypto market data (prices / 1h volume %) scraped on
# Wed Jul 26 03:43:09 UTC 2023.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value to a cell, forcing text storage (leading
# apostrophe, same as typing it in Excel) when the text would otherwise
# be auto-converted to a number and lose its original formatting
# (e.g. "0.9990" -> 0.999, "238.21" -> 238.21000000000001).
function Set-TextCell($row, $col, $text) {
    if ($text -match '^\s*[+-]?(\d+\.?\d*|\.\d+)\s*$') {
        $ws.Cells.Item($row, $col).Value = "'" + $text
    } else {
        $ws.Cells.Item($row, $col).Value = $text
    }
}

# Row 2
Set-TextCell 2 4 '29.221.38'
Set-TextCell 2 5 '  +0.35%  '
# Row 3
Set-TextCell 3 4 '1.858.72'
Set-TextCell 3 5 '  +0.36%  '
# Row 4
Set-TextCell 4 4 '0.9989'
Set-TextCell 4 5 '  -0.14%  '
# Row 5
Set-TextCell 5 4 '0.7106'
Set-TextCell 5 5 '  +2.07%  '
# Row 6
Set-TextCell 6 4 '238.21'
Set-TextCell 6 5 '  -0.14%  '
# Row 7
Set-TextCell 7 4 '0.9994'
Set-TextCell 7 5 '  -0.11%  '
# Row 8
Set-TextCell 8 4 '0.08011'
Set-TextCell 8 5 '  +4.83%  '
# Row 9
Set-TextCell 9 5 '  +0.10%  '
# Row 10
Set-TextCell 10 4 '23.58'
Set-TextCell 10 5 '  +0.83%  '
# Row 11
Set-TextCell 11 4 '0.08189'
Set-TextCell 11 5 '  +0.69%  '
# Row 12
Set-TextCell 12 4 '1.846.42'
Set-TextCell 12 5 '  -0.75%  '
# Row 13
Set-TextCell 13 4 '5.175'
Set-TextCell 13 5 '  -0.98%  '
# Row 14
Set-TextCell 14 4 '0.7045'
Set-TextCell 14 5 '  -3.09%  '
# Row 15
Set-TextCell 15 4 '89.71'
# Row 16
Set-TextCell 16 4 '29.193.01'
Set-TextCell 16 5 '  +0.26%  '
# Row 17
Set-TextCell 17 4 '5.842'
Set-TextCell 17 5 '  +1.02%  '
# Row 18
Set-TextCell 18 4 '0.000007879'
Set-TextCell 18 5 '  +1.81%  '
# Row 19
Set-TextCell 19 5 '  +0.68%  '
# Row 20
Set-TextCell 20 4 '238.17'
Set-TextCell 20 5 '  +0.65%  '
# Row 21
Set-TextCell 21 4 '0.9988'
Set-TextCell 21 5 '  -0.16%  '
# Row 22
Set-TextCell 22 2 'BinanceUSD'
Set-TextCell 22 3 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextCell 22 4 '0.9994'
Set-TextCell 22 5 '  -0.10%  '
# Row 23
Set-TextCell 23 2 'WrappedliquidstakedEther2.0'
Set-TextCell 23 3 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextCell 23 4 '2.069.07'
Set-TextCell 23 5 '  -1.33%  '
# Row 24
Set-TextCell 24 4 '7.435'
Set-TextCell 24 5 '  -2.38%  '
# Row 25
Set-TextCell 25 4 '162.30'
Set-TextCell 25 5 '  +0.54%  '
# Row 26
Set-TextCell 26 4 '8.934'
Set-TextCell 26 5 '  -0.61%  '
# Row 27
Set-TextCell 27 4 '0.1443'
Set-TextCell 27 5 '  -0.27%  '
# Row 28
Set-TextCell 28 5 '  -0.01%  '
# Row 29
Set-TextCell 29 4 '1.932'
Set-TextCell 29 5 '  -2.58%  '
# Row 30
Set-TextCell 30 5 '  +1.53%  '
# Row 31
Set-TextCell 31 4 '1.482'
Set-TextCell 31 5 '  -0.60%  '
# Row 32
Set-TextCell 32 4 '4.373'
Set-TextCell 32 5 '  -2.30%  '
# Row 33
Set-TextCell 33 4 '4.013'
Set-TextCell 33 5 '  -0.10%  '
# Row 34
Set-TextCell 34 5 '  -0.25%  '
# Row 35
Set-TextCell 35 5 '  -2.22%  '
# Row 36
Set-TextCell 36 4 '0.7098'
Set-TextCell 36 5 '  +1.26%  '
# Row 37
Set-TextCell 37 4 '0.9999'
Set-TextCell 37 5 '  -2.73%  '
# Row 38
Set-TextCell 38 4 '2.670'
Set-TextCell 38 5 '  +0.66%  '
# Row 39
Set-TextCell 39 4 '0.01858'
Set-TextCell 39 5 '  +0.18%  '
# Row 40
Set-TextCell 40 4 '2.725'
Set-TextCell 40 5 '  +1.73%  '
# Row 41
Set-TextCell 41 4 '0.9297'
Set-TextCell 41 5 '  +0.30%  '
# Row 42
Set-TextCell 42 4 '1.130.60'
Set-TextCell 42 5 '  +4.53%  '
# Row 43
Set-TextCell 43 5 '  -0.02%  '
# Row 44
Set-TextCell 44 4 '5.853'
Set-TextCell 44 5 '  -2.18%  '
# Row 45
Set-TextCell 45 4 '70.13'
Set-TextCell 45 5 '  -0.31%  '
# Row 46
Set-TextCell 46 4 '0.9990'
Set-TextCell 46 5 '  -0.16%  '
# Row 47
Set-TextCell 47 4 '102.94'
Set-TextCell 47 5 '  -0.24%  '
# Row 48
Set-TextCell 48 4 '0.5340'
Set-TextCell 48 5 '  -4.28%  '
# Row 49
Set-TextCell 49 4 '1.764'
Set-TextCell 49 5 '  -0.74%  '
# Row 50
Set-TextCell 50 2 'EnergySwap'
Set-TextCell 50 3 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 50 4 '9.161'
Set-TextCell 50 5 '  -0.37%  '
# Row 51
Set-TextCell 51 2 'RocketPoolETH'
Set-TextCell 51 3 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextCell 51 4 '1.967.42'
Set-TextCell 51 5 '  -1.32%  '
